$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 4 to row 5 first, to keep consistent styling
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the new journal entry row (row 5) as raw serial values, matching the
# underlying number-formatted cells (date serial, decimal hours, text)
$ws.Range("A5").Value = 43882
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "J'ai travaillé à la maison pour le MLD"

# Update selection to match the new active cell
$ws.Range("C5").Select()
